$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "'1-Car Garage'",
    "'2-Car Garage'",
    "'3-Car Garage'",
    "'Air Conditioning'",
    "'Attached Garage'",
    "'Balcony'",
    "'Basement'",
    "'Bathtub w/jets'",
    "'Bay Window'",
    "'Bonus Room'",
    "'Breakfast Nook'",
    "'Built-In BBQ'",
    "'Built-In Bookshelves'",
    "'Cable Ready'",
    "'Carpet'",
    "'Ceiling Fan'",
    "'Ceiling Fans'",
    "'Central Air/Heat'",
    "'Coffee System'",
    "'Covered Patio'",
    "'Crown Molding'",
    "'Deck'",
    "'Den'",
    "'Dining Room'",
    "'Dishwasher'",
    "'Disposal'",
    "'Dock'",
    "'Double Pane Windows'",
    "'Double Vanities'",
    "'Dual-Vanity Sinks'",
    "'Eat-in Kitchen'",
    "'Enclosed Porch / Sunroom'",
    "'Extended Driveway'",
    "'Family Room'",
    "'Fenced Pool'",
    "'Fenced Yard'",
    "'Fireplace'",
    "'Framed Mirrors'",
    "'Freezer'",
    "'Furnished'",
    "'Furnished Units Available'",
    "'Garage'",
    "'Garden'",
    "'Garden Tub'",
    "'Gated'",
    "'Granite Countertops'",
    "'Great Room'",
    "'Greenhouse'",
    "'Grill'",
    "'Handrails'",
    "'Hardwood Flooring'",
    "'Hardwood Floors'",
    "'Heating'",
    "'High Ceilings'",
    "'High Speed Internet Access'",
    "'Hookup: In-Unit, 1st Floor'",
    "'Hookup: In-Unit, 2nd Floor'",
    "'Ice Maker'",
    "'In Unit Washer & Dryer'",
    "'Instant Hot Water'",
    "'Intercom'",
    "'Island Kitchen'",
    "'Kitchen'",
    "'Kitchen Island'",
    "'Laminate Flooring'",
    "'Large Backyard'",
    "'Large Bedrooms'",
    "'Laundry Facilities'",
    "'Laundry Room'",
    "'Lawn'",
    "'Linen Closet'",
    "'Loft'",
    "'Loft Layout'",
    "'Microwave'",
    "'Mud Room'",
    "'Natural Light / Sky Lights'",
    "'Office'",
    "'Open Floor Plan'",
    "'Oven'",
    "'Oversized Bathtub'",
    "'Pantry'",
    "'Patio'",
    "'Plantation Shutters'",
    "'Pool'",
    "'Porch'",
    "'Quartz Countertops'",
    "'Range'",
    "'Rear-Loaded Driveway'",
    "'Recreation Room'",
    "'Refrigerator'",
    "'Satellite TV'",
    "'Screened Porch'",
    "'Security System'",
    "'Skylight'",
    "'Smart Home'",
    "'Smoke Free'",
    "'Sprinkler System'",
    "'Stainless Steel Appliances'",
    "'Storage Space'",
    "'Storage Units'",
    "'Sunroom'",
    "'Surround Sound'",
    "'Tile Flooring'",
    "'Tile Floors'",
    "'Trash Compactor'",
    "'Travertine Flooring'",
    "'Tub/Shower'",
    "'Vaulted Ceiling'",
    "'Views'",
    "'Vinyl Flooring'",
    "'Vinyl Plank Flooring'",
    "'Walk-In Closet'",
    "'Walk-In Closets'",
    "'Walk-In Pantry'",
    "'Walk-In Shower'",
    "'Wall Oven'",
    "'Warming Drawer'",
    "'Washer/Dryer Hookup'",
    "'Wet Bar'",
    "'Wheelchair Accessible (Rooms)'",
    "'Wi-Fi'",
    "'Window Coverings'",
    "'Wood Fireplace'",
    "'Wood-Look Blinds'",
    "'Yard'"
)

# Write each value as a formula that evaluates to the literal text (so a leading
# apostrophe in the text is not misinterpreted by Excel as a "text quote prefix").
for ($i = 0; $i -lt $values.Length; $i++) {
    $escaped = $values[$i].Replace("""", """""")
    $ws.Cells.Item($i + 1, 1).Formula = "=""" + $escaped + """"
}

# Convert the formulas to static values in one shot via copy / paste-special (values only),
# which stores them as plain shared-string cells, just like the target file.
$rng = $ws.Range("A1:A" + $values.Length)
$rng.Copy()
$rng.PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Clear the now-unused trailing rows left over from the longer original list (the sheet
# previously held 131 rows; the refreshed amenity list only needs $($values.Length)).
$lastOldRow = 131
if ($lastOldRow -gt $values.Length) {
    $ws.Range("A" + ($values.Length + 1) + ":A" + $lastOldRow).ClearContents()
}

